# Append the 01/19/2026 daily profit-tracking row to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first free row right after the current data block (row 55 -> 56).
$lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row()
$newRow = $lastRow + 1

# Column A holds the date as literal text (matches the rest of the sheet,
# which stores dates as plain "MM/DD/YYYY" strings, not Excel date serials).
# Force a text number format before assigning so the engine doesn't
# auto-convert the string into a date value, then clear the format again so
# the new cell ends up without any explicit style, same as its neighbours.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "01/19/2026"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value  = 12275.14
$ws.Cells.Item($newRow, 3).Value  = 0.2343547161733963
$ws.Cells.Item($newRow, 4).Value  = 0.7656452838266037
$ws.Cells.Item($newRow, 5).Value  = -165.03
$ws.Cells.Item($newRow, 6).Value  = -24.09
$ws.Cells.Item($newRow, 7).Value  = -21343.35
$ws.Cells.Item($newRow, 8).Value  = -69.43000000000001
$ws.Cells.Item($newRow, 9).Value  = -326.11
$ws.Cells.Item($newRow, 10).Value = -10.18
$ws.Cells.Item($newRow, 11).Value = -21669.46
$ws.Cells.Item($newRow, 12).Value = -63.84
